$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B and C values for rows 1-29 (row 1 is header)
$colB = @(
    "Folder Path",
    "Website ID 56+\+SSA-3288\",
    "Website ID 56+\+SSA-827\",
    "Website ID 56+\+SSA-1696\",
    "Website ID 56+\+SSA-8000\",
    "Website ID 56+\+Fillable SSA-8000\",
    '"Website ID 56"+"\"+"SSA-8001\"',
    '"Website ID 56"+"\"+"SSA-3369\"',
    '"Website ID 56"+"\"+"SSA-821\"',
    '"Website ID 56"+"\"+"HA-1151\"',
    '"Website ID 56"+"\"+"HA-1152\"',
    '"Website ID 56"+"\"+"SSA-787\"',
    '"Website ID 56"+"\"+"SSA-1699\"',
    '"Website ID 56"+"\"+"paper form\"',
    '"Website ID 56"+"\"+"paper form\"',
    '"Website ID 56"+"\"+"SSA-3373\"',
    '"Website ID 56"+"\"+"SSA-4814\"',
    '"Website ID 56"+"\"+"Listing of Impairments\"',
    '"Website ID 56"+"\"+"Online Disability Appeal Application\"',
    '"Website ID 56"+"\"+"SSA-3820\"',
    '"Website ID 56"+"\"+"SSA-3375\"',
    '"Website ID 56"+"\"+"3376\"',
    '"Website ID 56"+"\"+"3377\"',
    '"Website ID 56"+"\"+"3378\"',
    '"Website ID 56"+"\"+"3379\"',
    '"Website ID 56"+"\"+"SSA-5665\"',
    '"Website ID 56"+"\"+"Listing of Impairments\"',
    '"Website ID 56"+"\"+"Medicaid Eligibility Income\"',
    '"Website ID 56"+"\"+"Compassionate Allowances\"'
)

$colC = @(
    "Folder Path 2",
    '"Website ID 56+"\"+"SSA-3288\"',
    '"Website ID 56"+"\"+"SSA-827\"',
    '"Website ID 56"+"\"+"SSA-1696\"',
    '"Website ID 56"+"\"+"SSA-8000\"',
    '"Website ID 56"+"\"+"Fillable SSA-8000\"',
    '"Website ID 56"+"\"+"SSA-8001\"',
    '"Website ID 56"+"\"+"SSA-3369\"',
    '"Website ID 56"+"\"+"SSA-821\"',
    '"Website ID 56"+"\"+"HA-1151\"',
    '"Website ID 56"+"\"+"HA-1152\"',
    '"Website ID 56"+"\"+"SSA-787\"',
    '"Website ID 56"+"\"+"SSA-1699\"',
    '"Website ID 56"+"\"+"paper form\"',
    '"Website ID 56"+"\"+"paper form\"',
    '"Website ID 56"+"\"+"SSA-3373\"',
    '"Website ID 56"+"\"+"SSA-4814\"',
    '"Website ID 56"+"\"+"Listing of Impairments\"',
    '"Website ID 56"+"\"+"Online Disability Appeal Application\"',
    '"Website ID 56"+"\"+"SSA-3820\"',
    '"Website ID 56"+"\"+"SSA-3375\"',
    '"Website ID 56"+"\"+"3376\"',
    '"Website ID 56"+"\"+"3377\"',
    '"Website ID 56"+"\"+"3378\"',
    '"Website ID 56"+"\"+"3379\"',
    '"Website ID 56"+"\"+"SSA-5665\"',
    '"Website ID 56"+"\"+"Listing of Impairments\"',
    '"Website ID 56"+"\"+"Medicaid Eligibility Income\"',
    '"Website ID 56"+"\"+"Compassionate Allowances\"'
)

for ($i = 0; $i -lt 29; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $colB[$i]
    $ws.Cells.Item($row, 3).Value = $colC[$i]
}

$ws.Columns.Item(2).ColumnWidth = 56.33
$ws.Columns.Item(3).ColumnWidth = 52.5

$ws.Range("B9").Select()
